$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (C) column date serial from 45174 to 45175 (2023-09-05 -> 2023-09-06)
# for all rows currently holding the old value (rows 2-10).
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45174) {
        $cell.Value = 45175
    }
}
